$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 915763.7357711862
$ws.Range("C3").Value = 171037.82777760012
$ws.Range("C4").Value = 86921.60302080003
$ws.Range("C5").Value = 657804.3049727976
$ws.Range("C6").Value = 158705.82012960006
$ws.Range("C7").Value = 230822.24114880018
$ws.Range("C8").Value = 272825.2078272002
$ws.Range("C9").Value = 32494.692268800005
$ws.Range("C10").Value = 220747.3237008001
$ws.Range("C11").Value = 168.450696
